$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting of the last existing data row (27) down into the new row (28),
# matching each column's existing look. Column B keeps the default (unstyled)
# look, same as in the source file, so it is intentionally left untouched.
$ws.Range("A27").Copy()
$ws.Range("A28").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C27").Copy()
$ws.Range("C28").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D27").Copy()
$ws.Range("D28").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new timeline entry
$ws.Cells.Item(28, 1).Value = 25
$ws.Cells.Item(28, 2).Value = "29/3/2024"
$ws.Cells.Item(28, 3).Value = 5.4
$ws.Cells.Item(28, 4).Value = "Added batch processing of data + product recommendations"

# Restore the view state (scroll position / selection) as in the saved workbook
$ws.Activate() | Out-Null
$ws.Range("C30").Select() | Out-Null
